$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data block of interest spans rows 285..404 (120 rows), columns A..R (18 cols).
# The edit "shifts" the record fields D,J,K,L,M,N,O,P,Q down by one row within this
# block: each row (from 286 downward) receives the shiftable fields that used to
# belong to the row immediately above it. Row 285 receives brand-new field values
# (D and J), and a new row 405 is appended holding the fields that "fell off the
# bottom" (i.e. the original row 404 values).

$startRow = 285
$endRow = 404
$numOldRows = $endRow - $startRow + 1      # 120
$numNewRows = $numOldRows + 1              # 121

# Read the existing block as a 2-D array (Excel COM safearrays coming back from a
# range read are 1-based: $old[1..numOldRows, 1..18]).
$srcRange = $ws.Range("A" + $startRow + ":R" + $endRow)
$old = $srcRange.Value2

# Columns (1-based within the A..R block) that shift down by one row:
# D=4, J=10, K=11, L=12, M=13, N=14, O=15, P=16, Q=17
$shiftCols = @(4,10,11,12,13,14,15,16,17)

# Build the replacement block. Arrays created with New-Object are 0-based, and that
# is also what the Range.Value2 *setter* expects (array index [0,0] lands on the
# top-left cell of the destination range).
$new = New-Object 'object[,]' $numNewRows, 18

# --- New row for sheet row 285 (array row 0): keep everything from the old row 285,
#     except D and J which take brand-new values.
for ($c = 1; $c -le 18; $c++) {
    $new[0, $c - 1] = $old[1, $c]
}
$new[0, 3] = 44875   # column D (0-based index 3)
$new[0, 9] = 95      # column J (0-based index 9)

# --- Rows for sheet rows 286..404 (array rows 1..numOldRows-1): non-shifted columns
#     keep their own row's original values; shifted columns take the values that
#     belonged to the previous row.
for ($i = 1; $i -le ($numOldRows - 1); $i++) {
    $ownOldRow = $i + 1   # 1-based old-array row holding this sheet row's own data
    for ($c = 1; $c -le 18; $c++) {
        $new[$i, $c - 1] = $old[$ownOldRow, $c]
    }
    foreach ($c in $shiftCols) {
        $new[$i, $c - 1] = $old[$i, $c]   # previous row's (shifted) values
    }
}

# --- Brand-new last row for sheet row 405 (array row numNewRows-1): every column
#     (shifted or not) takes the values that "overflowed" from the old row 404.
$lastIdx = $numNewRows - 1
for ($c = 1; $c -le 18; $c++) {
    $new[$lastIdx, $c - 1] = $old[$numOldRows, $c]
}

$newEndRow = $endRow + 1
$destRange = $ws.Range("A" + $startRow + ":R" + $newEndRow)
$destRange.Value2 = $new

# The new D405 cell was created from scratch and therefore has no number format;
# give it the same (date) number format used by the rest of column D.
$ws.Range("D" + $newEndRow).NumberFormat = $ws.Range("D" + $endRow).NumberFormat
